$d = $word.ActiveDocument

function FindRange($searchText) {
    $r = $d.Content
    $found = $r.Find.Execute($searchText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $found) {
        throw "Text not found: $searchText"
    }
    return $r
}

# ---------------------------------------------------------------------------
# Change 1 (Problem 2 / KLT paragraph): split the trailing sentence about
# de-correlating the input vector into its own new paragraph, and insert a
# new sentence about the final misadjustment before it.
# ---------------------------------------------------------------------------
$oldTail = " By de-correlating the input vector, a faster convergence is possible for the same normalized step size. The drawback for KLT is required knowledge of the eigenvectors for linearly mapping the input vector to the decorrelated transformed domain, and the eigenvalues for proper scaling of the step size for each filter tap."

$r1 = FindRange($oldTail)
$splitPoint = $r1.Start

# Remove the whole old run's text first (it sits at the very end of its
# paragraph, so deleting it cleanly drops the run without merging anything).
$r1.Delete()

# Re-insert the pieces as brand new runs, from the split point onward.
$ins = $d.Range($splitPoint, $splitPoint)
$ins.InsertAfter(" ")
$ins.Collapse(0)
$ins.InsertAfter("Final misadjustment was also comparable to the LMS. ")
$ins.Collapse(0)
$ins.InsertAfter("`r")
$ins.Collapse(0)
$ins.InsertAfter("By de-correlating the input vector, a faster convergence is possible for the same normalized step size. The drawback for KLT is required knowledge of the eigenvectors for linearly mapping the input vector to the decorrelated transformed domain, and the eigenvalues for proper scaling of the step size for each filter tap.")

# ---------------------------------------------------------------------------
# Change 2 (Problem 6 / BIC paragraph): "low" -> "high", and move the
# "_GoBack" bookmark to sit right after the new word "high". Word only keeps
# a single "_GoBack" bookmark, so adding it here automatically removes the
# one that used to sit at the end of the document (Change 3).
# ---------------------------------------------------------------------------
$r2 = FindRange("SNR should be low")
$lowStart = $r2.End - 3
$lowEnd = $r2.End

$lowRange = $d.Range($lowStart, $lowEnd)
$lowRange.Delete()

$insHigh = $d.Range($lowStart, $lowStart)
$insHigh.InsertAfter("high")
$highEnd = $lowStart + 4

# Force a run split before "high" using a temporary bookmark, then remove it.
$beforeHigh = $d.Range($lowStart, $lowStart)
$tempBookmark = $d.Bookmarks.Add("ZZZTempSplit", $beforeHigh)

# Force a run split after "high" by placing the real "_GoBack" bookmark there.
$afterHigh = $d.Range($highEnd, $highEnd)
$d.Bookmarks.Add("_GoBack", $afterHigh) | Out-Null

$tempBookmark = $d.Bookmarks.Item("ZZZTempSplit")
$tempBookmark.Delete()
